$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J.
# Copy an existing header cell (H1) onto I1/J1 first so the destination
# picks up the exact same cell style (bold/border/centered) that the rest
# of the header row uses, then overwrite the copied text with the real
# header labels.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# New data values for columns I and J (rows 2-11)
$values = @(
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(7, 7),
    @(3, 3),
    @(10, 10),
    @(8, 8),
    @(8, 8),
    @(7, 7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
